$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: metadata type
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "iaest-measure:estudios-en-curso"
$ws.Range("F2").Value = "sdmx-dimension:refArea"

# Row 3: medida/dim classification
$ws.Range("C3").Value = "dim"
$ws.Range("E3").Value = "medida"

# Row 4: data type / URI
$ws.Range("C4").Value = "URI-Municipio"
$ws.Range("E4").Value = "xsd:int"
$ws.Range("F4").Value = "URI-Comunidad"

# Row 5 (E5, F5) no longer present - remove the row entirely
$ws.Range("A5:I5").EntireRow.Delete() | Out-Null
